$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update sprint review table (rows 16-20) ---
# Row 16: Grid Layour erstellen - estimate 5 -> 10, real effort -> 10, passed
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 10
$ws.Range("F16").Value = 1

# Row 17: Itemliste rechts im Layout hinzufuegen - real effort -> 2, passed, remark
$ws.Range("D17").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("H17").Value = "Pair programming"

# Row 18: Item Icons erstellen - real effort 0.5 -> 1, passed
$ws.Range("D18").Value = 1
$ws.Range("F18").Value = 1

# Row 19: Mock Ups fuer die einzelnen Seite erstellen - real effort -> 4, passed, remark
$ws.Range("D19").Value = 4
$ws.Range("F19").Value = 1
$ws.Range("H19").Value = "Pair programming"

# Row 20: Ordnerstruktur anlegen - real effort -> 7, passed, remark
$ws.Range("D20").Value = 7
$ws.Range("F20").Value = 1
$ws.Range("H20").Value = "Pair programming"

# --- Column B widened to fit longer text ---
$ws.Columns("B").ColumnWidth = 41.88671875

# --- Page setup: landscape, scaled to 90%, metric-rounded margins ---
$ps = $ws.PageSetup
$ps.Orientation = 2
$ps.Zoom = 90
$ps.LeftMargin = 51.0236220472441
$ps.RightMargin = 51.0236220472441
$ps.TopMargin = 56.69291338582678
$ps.BottomMargin = 56.69291338582678
$ps.HeaderMargin = 22.677165354330707
$ps.FooterMargin = 22.677165354330707

# Manual column break before column I
$ws.Columns("I").PageBreak = 1

# --- View / selection state ---
$excel.ActiveWindow.Zoom = 60
$ws.Range("E18").Select()
